# Automatische test-sync: 2025-07-27 16:23:50
# Append a new log entry (row 4) to the "Logs" sheet, mirroring the
# existing rows 2/3, extend the conditional-formatting ranges to cover
# the new row, and bump the "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$logs.Range("A4").Value = "Kun jij dit even regelen?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D4").Value = "Overig"
$logs.Range("E4").Value = "Geachte klant,`nBedankt voor uw e-mail. Om u beter van dienst te kunnen zijn, zou u wat meer informatie kunnen verstrekken over wat u precies wilt regelen?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F4").Value = "2025-07-27 16:23:23"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Nee"
$logs.Range("I4").Value = "Ja"
$logs.Range("J4").Value = "Ja"

# Extend the conditional formatting sqref for each column from row 2:3 to 2:4
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $fcs = $logs.Range($col + "2").FormatConditions
    $count = $fcs.Count
    for ($i = 1; $i -le $count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($logs.Range($col + "2:" + $col + "4"))
    }
}

$dashboard.Range("B2").Value = 3
